$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Remove the row describing "Cadastro de valores de serviços prestados." (row 7)
$ws.Rows.Item(7).Delete()

# Sort the remaining data rows (2:8) ascending by column B ("Caracteristicas")
$rng = $ws.Range("A1:E8")
$rng.Sort($ws.Range("B1"), 1)

# Turn on AutoFilter for the header row
$ws.Range("A1:E1").AutoFilter() | Out-Null

# Record the hidden sheet-scoped _FilterDatabase defined name that Excel
# creates when an AutoFilter is applied
$nm = $ws.Names.Add("_xlnm._FilterDatabase", "=Plan1!`$A`$1:`$E`$1")
$nm.Visible = $false

# Update the active selection to match the final workbook state
$ws.Range("A13").Select() | Out-Null
